$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.440.66"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.626.04"
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.87"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.51"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.207"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.34"
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.57"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.198.49"
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "599.54"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.01"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.543.95"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.604.46"
$ws.Range("E18").Value = "  +2.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.04"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.15"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.20"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.99"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E26").Value = "  -3.35%  "
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.83"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("E30").Value = "  +6.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.26"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.28"
$ws.Range("E32").Value = "  -3.22%  "
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.48"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("E35").Value = "  +4.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.918.97"
$ws.Range("E36").Value = "  +4.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "538.85"
$ws.Range("E37").Value = "  +9.38%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.98"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.390"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("E42").Value = "  -3.36%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.134"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0463"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.88"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.39"
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000249"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("E51").Value = "  +1.59%  "
